# "add two insert cases"
#
# Appends two new test-case rows (insert_018, insert_019) to the bottom of
# Sheet1's test matrix, following the exact same column layout used by every
# other row in the sheet (TestID, Testable, Title, Component, Sub_component,
# Table_schema_ref, Table_value_ref, Effected_rows, Query_sql1,
# Expected_result1, ..., Validation_type), and moves the saved selection to
# the author's final cursor position.
#
# NOTE: this engine's PowerShell function calls only bind POSITIONAL
# arguments reliably, so Add-TestCaseRow below takes its values in a fixed
# order instead of via -Name value switches.

function Add-TestCaseRow($Row, $TemplateRow, $TestID, $Testable, $Title, $Component, $SchemaRef, $ValueRef, $EffectedRows, $QuerySql1, $ExpectedResult1, $ValidationType) {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.Worksheets.Item(1)

    # Columns that actually carry content on the template row: A-J, M, N, O.
    # (K/L are intentionally skipped -- the template row has nothing in them.)
    $templateCols = @(1,2,3,4,5,6,7,8,9,10,13,14,15)

    # Clone the formatting of the template row (per populated column) so the
    # new row renders identically to every other case row.
    foreach ($c in $templateCols) {
        $ws.Cells.Item($TemplateRow, $c).Copy()
        $ws.Cells.Item($Row, $c).PasteSpecial(-4122)  # xlPasteFormats
    }
    $excel.CutCopyMode = $false

    $ws.Cells.Item($Row, 1).Value  = $TestID
    $ws.Cells.Item($Row, 2).Value  = $Testable
    $ws.Cells.Item($Row, 3).Value  = $Title
    $ws.Cells.Item($Row, 4).Value  = $Component
    $ws.Cells.Item($Row, 6).Value  = $SchemaRef
    $ws.Cells.Item($Row, 7).Value  = $ValueRef
    $ws.Cells.Item($Row, 8).Value  = $EffectedRows
    $ws.Cells.Item($Row, 9).Value  = $QuerySql1
    $ws.Cells.Item($Row, 10).Value = $ExpectedResult1
    $ws.Cells.Item($Row, 15).Value = $ValidationType
}

Add-TestCaseRow 117 116 "insert_018" "y" "指定非连续字段插入数据1" "insert" `
    "schema7" "insert_value15" "2" 'select * from $schema7' `
    "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/insert/expectedresult/insert_018.csv" `
    "csv_containsAll"

Add-TestCaseRow 118 117 "insert_019" "y" "指定非连续字段插入数据2" "insert" `
    "schema17" "insert_value16" "3" 'select * from $schema17' `
    "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/insert/expectedresult/insert_019.csv" `
    "csv_containsAll"

# Leave the selection where the author last left it when saving.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("M116").Select()
